# Daily attendance processing - reorders the "Recorded By" (column G)
# actor list so that any "System"/"system" entry that currently leads the
# comma-separated list is moved to the end (the whole list is reversed
# whenever the first token is "System", case-insensitively).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -ne $null -and $val -ne "") {
        $parts = $val.Split(", ")

        if ($parts.Count -gt 1 -and $parts[0].ToLower() -eq "system") {
            $n = $parts.Count
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $newVal = $reversed -join ", "
            $cell.Value = $newVal
        }
    }
}
